# Update the MSME indicator figures for Egypt, Arab Rep. with more precise
# (two decimal place) values. These cells hold text values (e.g. "33.5"),
# so a leading apostrophe is used to force Excel to keep storing them as
# text rather than auto-converting to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Enterprises density (per 1000 people): Micro / SMEs / MSMEs
$ws.Range("B13").Value = "'33.53"
$ws.Range("C13").Value = "'0.08"
$ws.Range("D13").Value = "'33.61"

# Employment (% of total): Micro / SMEs / MSMEs
$ws.Range("B14").Value = "'74.29"
$ws.Range("C14").Value = "'20.53"
$ws.Range("D14").Value = "'94.82"

# Enterprises (% of total): Micro / SMEs / MSMEs
$ws.Range("B16").Value = "'99.73"
$ws.Range("C16").Value = "'0.25"
$ws.Range("D16").Value = "'99.98"
